# Add season record columns (Wins, Losses, Ties) to the BOS_2015 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row (header row 1 + data rows below).
$lastRow = $ws.UsedRange.Rows.Count

# New columns: AD (30) = Wins, AE (31) = Losses, AF (32) = Ties
$winsCol   = 30
$lossesCol = 31
$tiesCol   = 32

# Copy the header formatting from the existing last header cell (AC1) onto
# the new header cells so they match the bold/centered/bordered style used
# by the rest of row 1.
$ws.Cells.Item(1, 29).Copy()
$ws.Cells.Item(1, $winsCol).PasteSpecial(-4122)
$ws.Cells.Item(1, $lossesCol).PasteSpecial(-4122)
$ws.Cells.Item(1, $tiesCol).PasteSpecial(-4122)

$ws.Cells.Item(1, $winsCol).Value   = "Wins"
$ws.Cells.Item(1, $lossesCol).Value = "Losses"
$ws.Cells.Item(1, $tiesCol).Value   = "Ties"

# Season record for every player row: 78 wins, 84 losses, 0 ties.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $winsCol).Value   = 78
    $ws.Cells.Item($r, $lossesCol).Value = 84
    $ws.Cells.Item($r, $tiesCol).Value   = 0
}
